$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '26.254.07'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '1.592.07'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'212.51"
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('E6').Value = '  -0.85%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.87%  '
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D10').Value = "'18.98"
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('D11').Value = "'0.0851"
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '1.815.21'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '1.605.83'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('E15').Value = '  -3.09%  '
$ws.Range('D16').Value = "'63.86"
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').Value = '26.238.12'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = "'214.68"
$ws.Range('D20').Value = "'7.33"
$ws.Range('E20').Value = '  -2.62%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').Value = "'4.30"
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').Value = "'9.05"
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').Value = '  -2.46%  '
$ws.Range('D25').Value = "'144.47"
$ws.Range('E25').Value = '  +0.34%  '
$ws.Range('D27').Value = "'6.96"
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('D29').Value = "'15.12"
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('E30').Value = '  -2.80%  '
$ws.Range('D31').Value = "'1.16"
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').Value = "'3.20"
$ws.Range('E32').Value = '  -0.80%  '
$ws.Range('D33').Value = '1.410.10'
$ws.Range('E33').Value = '  +5.68%  '
$ws.Range('D34').Value = "'2.96"
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('D35').Value = "'2.42"
$ws.Range('E35').Value = '  -0.88%  '
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('D37').Value = "'0.582"
$ws.Range('E37').Value = '  -3.51%  '
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('D39').Value = "'0.823"
$ws.Range('E39').Value = '  +0.52%  '
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = "'0.985"
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = "'0.764"
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '1.727.71'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').Value = "'60.90"
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('D47').Value = "'86.47"
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').Value = "'0.0953"
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('E51').Value = '  +0.02%  '
